{"js": "// Update the date line and all of the two-digit-by-two-digit multiplication\n// prompts in the practice sheet to the new values from the latest generator\n// run.\nconst replacements = [\n  [\"2024-05-22 Wednesday\", \"2024-05-23 Thursday\"],\n  [\"84\u00d712=\", \"66\u00d742=\"],\n  [\"14\u00d720=\", \"87\u00d797=\"],\n  [\"57\u00d756=\", \"84\u00d734=\"],\n  [\"89\u00d728=\", \"26\u00d743=\"],\n  [\"43\u00d743=\", \"59\u00d733=\"],\n  [\"64\u00d788=\", \"27\u00d796=\"],\n  [\"53\u00d750=\", \"95\u00d731=\"],\n  [\"70\u00d735=\", \"90\u00d711=\"],\n  [\"36\u00d777=\", \"88\u00d713=\"],\n  [\"38\u00d737=\", \"68\u00d741=\"],\n  [\"52\u00d797=\", \"81\u00d762=\"],\n  [\"35\u00d765=\", \"80\u00d784=\"],\n  [\"43\u00d770=\", \"82\u00d799=\"],\n  [\"58\u00d799=\", \"94\u00d755=\"],\n  [\"35\u00d750=\", \"75\u00d764=\"],\n  [\"40\u00d741=\", \"16\u00d778=\"],\n  [\"60\u00d768=\", \"85\u00d712=\"],\n  [\"46\u00d777=\", \"64\u00d797=\"],\n  [\"49\u00d723=\", \"93\u00d747=\"],\n  [\"79\u00d711=\", \"89\u00d744=\"],\n  [\"74\u00d736=\", \"74\u00d794=\"],\n  [\"45\u00d738=\", \"32\u00d790=\"],\n  [\"48\u00d717=\", \"33\u00d788=\"],\n  [\"20\u00d718=\", \"55\u00d726=\"],\n  [\"42\u00d713=\", \"36\u00d782=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all of the two-digit-by-two-digit multiplication\n# prompts in the practice sheet to the new values from the latest generator\n# run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-22 Wednesday\", \"2024-05-23 Thursday\"),\n    @(\"84\u00d712=\", \"66\u00d742=\"),\n    @(\"14\u00d720=\", \"87\u00d797=\"),\n    @(\"57\u00d756=\", \"84\u00d734=\"),\n    @(\"89\u00d728=\", \"26\u00d743=\"),\n    @(\"43\u00d743=\", \"59\u00d733=\"),\n    @(\"64\u00d788=\", \"27\u00d796=\"),\n    @(\"53\u00d750=\", \"95\u00d731=\"),\n    @(\"70\u00d735=\", \"90\u00d711=\"),\n    @(\"36\u00d777=\", \"88\u00d713=\"),\n    @(\"38\u00d737=\", \"68\u00d741=\"),\n    @(\"52\u00d797=\", \"81\u00d762=\"),\n    @(\"35\u00d765=\", \"80\u00d784=\"),\n    @(\"43\u00d770=\", \"82\u00d799=\"),\n    @(\"58\u00d799=\", \"94\u00d755=\"),\n    @(\"35\u00d750=\", \"75\u00d764=\"),\n    @(\"40\u00d741=\", \"16\u00d778=\"),\n    @(\"60\u00d768=\", \"85\u00d712=\"),\n    @(\"46\u00d777=\", \"64\u00d797=\"),\n    @(\"49\u00d723=\", \"93\u00d747=\"),\n    @(\"79\u00d711=\", \"89\u00d744=\"),\n    @(\"74\u00d736=\", \"74\u00d794=\"),\n    @(\"45\u00d738=\", \"32\u00d790=\"),\n    @(\"48\u00d717=\", \"33\u00d788=\"),\n    @(\"20\u00d718=\", \"55\u00d726=\"),\n    @(\"42\u00d713=\", \"36\u00d782=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
